$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WIT")

$ws.Range("D8").Value = 7878800
$ws.Range("E8").Value = 7958800
$ws.Range("F8").Value = 7409900
$ws.Range("G8").Value = 6789600
$ws.Range("H8").Value = 6279500
$ws.Range("I8").Value = 5411700
$ws.Range("J8").Value = 4609100
$ws.Range("D9").Value = 5575400
$ws.Range("E9").Value = 5661700
$ws.Range("F9").Value = 5158200
$ws.Range("G9").Value = 4645800
$ws.Range("H9").Value = 4272800
$ws.Range("I9").Value = 3769200
$ws.Range("J9").Value = 3265000
$ws.Range("D10").Value = 2303400
$ws.Range("E10").Value = 2297100
$ws.Range("F10").Value = 2251700
$ws.Range("G10").Value = 2143900
$ws.Range("H10").Value = 2006800
$ws.Range("I10").Value = 1642500
$ws.Range("J10").Value = 1344100
$ws.Range("D14").Value = 9300
$ws.Range("E14").Value = -14800
$ws.Range("D15").Value = 55900
$ws.Range("E15").Value = 45900
$ws.Range("F15").Value = 25100
$ws.Range("D17").Value = 6659900
$ws.Range("E17").Value = 6601300
$ws.Range("F17").Value = 6009300
$ws.Range("G17").Value = 5409800
$ws.Range("H17").Value = 4987500
$ws.Range("I17").Value = 4399900
$ws.Range("J17").Value = 3742800
$ws.Range("D18").Value = 1218900
$ws.Range("E18").Value = 1357500
$ws.Range("F18").Value = 1400600
$ws.Range("G18").Value = 1379800
$ws.Range("H18").Value = 1292100
$ws.Range("I18").Value = 1011800
$ws.Range("J18").Value = 866300
$ws.Range("D20").Value = 312800
$ws.Range("E20").Value = 266000
$ws.Range("F20").Value = 281700
$ws.Range("G20").Value = 246200
$ws.Range("H20").Value = 181000
$ws.Range("I20").Value = 137200
$ws.Range("J20").Value = 94700
$ws.Range("D21").Value = 1837300
$ws.Range("E21").Value = 1957800
$ws.Range("F21").Value = 1898900
$ws.Range("G21").Value = 1811600
$ws.Range("H21").Value = 1633800
$ws.Range("I21").Value = 1305800
$ws.Range("J21").Value = 1107600
$ws.Range("D22").Value = 49900
$ws.Range("E22").Value = 27700
$ws.Range("F22").Value = 20400
$ws.Range("G22").Value = 11100
$ws.Range("H22").Value = 12600
$ws.Range("I22").Value = 12500
$ws.Range("J22").Value = 13500
$ws.Range("D23").Value = 1481800
$ws.Range("E23").Value = 1595700
$ws.Range("F23").Value = 1661900
$ws.Range("G23").Value = 1614900
$ws.Range("H23").Value = 1460500
$ws.Range("I23").Value = 1136500
$ws.Range("J23").Value = 947500
$ws.Range("D24").Value = 346400
$ws.Range("E24").Value = 364600
$ws.Range("F24").Value = 366800
$ws.Range("G24").Value = 356100
$ws.Range("H24").Value = 326800
$ws.Range("I24").Value = 244500
$ws.Range("J24").Value = 187300
$ws.Range("D26").Value = 1135400
$ws.Range("E26").Value = 1231200
$ws.Range("F26").Value = 1295100
$ws.Range("G26").Value = 1258900
$ws.Range("H26").Value = 1133700
$ws.Range("I26").Value = 892000
$ws.Range("J26").Value = 760100
$ws.Range("D27").Value = 1135400
$ws.Range("E27").Value = 1227600
$ws.Range("F27").Value = 1288000
$ws.Range("G27").Value = 1251200
$ws.Range("H27").Value = 1127400
$ws.Range("I27").Value = 887300
$ws.Range("J27").Value = 756600
$ws.Range("D29").Value = 22600
$ws.Range("I29").Value = 72300
$ws.Range("J29").Value = 49200
$ws.Range("D32").Value = -312800
$ws.Range("E32").Value = -266000
$ws.Range("F32").Value = -281700
$ws.Range("G32").Value = -246200
$ws.Range("H32").Value = -181000
$ws.Range("I32").Value = -137200
$ws.Range("J32").Value = -94700
$ws.Range("D33").Value = 1158000
$ws.Range("E33").Value = 1227600
$ws.Range("F33").Value = 1288000
$ws.Range("G33").Value = 1251200
$ws.Range("H33").Value = 1127400
$ws.Range("I33").Value = 959600
$ws.Range("J33").Value = 805900
$ws.Range("D35").Value = 1158000
$ws.Range("E35").Value = 1227600
$ws.Range("F35").Value = 1288000
$ws.Range("G35").Value = 1251200
$ws.Range("H35").Value = 1127400
$ws.Range("I35").Value = 959600
$ws.Range("J35").Value = 805900
$ws.Range("D41").Value = 649600
$ws.Range("E41").Value = 762200
$ws.Range("F41").Value = 1432200
$ws.Range("G41").Value = 2298300
$ws.Range("H41").Value = 1651300
$ws.Range("I41").Value = 516000
$ws.Range("J41").Value = 594900
$ws.Range("D42").Value = 3634700
$ws.Range("E42").Value = 4247900
$ws.Range("F42").Value = 6934300
$ws.Range("G42").Value = 1381900
$ws.Range("H42").Value = 1100000
$ws.Range("I42").Value = 1747500
$ws.Range("J42").Value = 1164500
$ws.Range("D43").Value = 2182200
$ws.Range("E43").Value = 2230800
$ws.Range("F43").Value = 3822500
$ws.Range("G43").Value = 2156500
$ws.Range("H43").Value = 2020700
$ws.Range("I43").Value = 1762800
$ws.Range("J43").Value = 1744700
$ws.Range("D44").Value = 75000
$ws.Range("E44").Value = 77500
$ws.Range("F44").Value = 99700
$ws.Range("G44").Value = 105300
$ws.Range("H44").Value = 68500
$ws.Range("I44").Value = 47200
$ws.Range("J44").Value = 154200
$ws.Range("D45").Value = 777400
$ws.Range("E45").Value = 474100
$ws.Range("F45").Value = 420000
$ws.Range("G45").Value = 369700
$ws.Range("H45").Value = 288200
$ws.Range("I45").Value = 371600
$ws.Range("J45").Value = 296300
$ws.Range("D46").Value = 7319000
$ws.Range("E46").Value = 7792500
$ws.Range("F46").Value = 7270800
$ws.Range("G46").Value = 6311700
$ws.Range("H46").Value = 5128800
$ws.Range("I46").Value = 4445000
$ws.Range("J46").Value = 3954600
$ws.Range("D47").Value = 514800
$ws.Range("E47").Value = 372800
$ws.Range("F47").Value = 303400
$ws.Range("G47").Value = 220900
$ws.Range("H47").Value = 186100
$ws.Range("I47").Value = 155200
$ws.Range("J47").Value = 231700
$ws.Range("D48").Value = 939300
$ws.Range("E48").Value = 1009200
$ws.Range("F48").Value = 939400
$ws.Range("G48").Value = 783800
$ws.Range("H48").Value = 744000
$ws.Range("I48").Value = 730600
$ws.Range("J48").Value = 853000
$ws.Range("D49").Value = 1962200
$ws.Range("E49").Value = 2049200
$ws.Range("F49").Value = 1703900
$ws.Range("G49").Value = 1099100
$ws.Range("H49").Value = 945100
$ws.Range("I49").Value = 816600
$ws.Range("J49").Value = 1043500
$ws.Range("D52").Value = 263500
$ws.Range("E52").Value = 250500
$ws.Range("F52").Value = 251700
$ws.Range("G52").Value = 261000
$ws.Range("H52").Value = 259500
$ws.Range("I52").Value = 211100
$ws.Range("J52").Value = 221700
$ws.Range("D54").Value = 10998900
$ws.Range("E54").Value = 11474200
$ws.Range("F54").Value = 10469100
$ws.Range("G54").Value = 8676500
$ws.Range("H54").Value = 7263300
$ws.Range("I54").Value = 6358500
$ws.Range("J54").Value = 6304600
$ws.Range("D57").Value = 352900
$ws.Range("E57").Value = 339100
$ws.Range("F57").Value = 1325000
$ws.Range("G57").Value = 272500
$ws.Range("H57").Value = 254700
$ws.Range("I57").Value = 223200
$ws.Range("J57").Value = 338800
$ws.Range("D58").Value = 1344600
$ws.Range("E58").Value = 1775700
$ws.Range("F58").Value = 3119300
$ws.Range("G58").Value = 957300
$ws.Range("H58").Value = 588300
$ws.Range("I58").Value = 910400
$ws.Range("J58").Value = 527000
$ws.Range("D59").Value = 1389800
$ws.Range("E59").Value = 1204400
$ws.Range("F59").Value = 1261700
$ws.Range("G59").Value = 1142000
$ws.Range("H59").Value = 1130200
$ws.Range("I59").Value = 959300
$ws.Range("J59").Value = 835900
$ws.Range("D60").Value = 3087300
$ws.Range("E60").Value = 3319200
$ws.Range("F60").Value = 3160400
$ws.Range("G60").Value = 2371900
$ws.Range("H60").Value = 1973200
$ws.Range("I60").Value = 2092900
$ws.Range("J60").Value = 1701700
$ws.Range("D61").Value = 654600
$ws.Range("E61").Value = 283600
$ws.Range("F61").Value = 251000
$ws.Range("G61").Value = 183700
$ws.Range("H61").Value = 157700
$ws.Range("I61").Value = 12300
$ws.Range("J61").Value = 325500
$ws.Range("D62").Value = 239000
$ws.Range("E62").Value = 313300
$ws.Range("F62").Value = 299300
$ws.Range("G62").Value = 197700
$ws.Range("H62").Value = 145400
$ws.Range("I62").Value = 132400
$ws.Range("J62").Value = 139400
$ws.Range("D66").Value = 4015600
$ws.Range("E66").Value = 3950600
$ws.Range("F66").Value = 3742700
$ws.Range("G66").Value = 2777100
$ws.Range("H66").Value = 2296300
$ws.Range("I66").Value = 2254600
$ws.Range("J66").Value = 2178900
$ws.Range("D72").Value = 6576600
$ws.Range("E72").Value = 7120200
$ws.Range("F72").Value = 6150300
$ws.Range("G72").Value = 5382700
$ws.Range("H72").Value = 4554200
$ws.Range("I72").Value = 3788300
$ws.Range("J72").Value = 3507000
$ws.Range("D76").Value = 6983300
$ws.Range("E76").Value = 7523600
$ws.Range("F76").Value = 6726400
$ws.Range("G76").Value = 5899400
$ws.Range("H76").Value = 4967000
$ws.Range("I76").Value = 4103900
$ws.Range("J76").Value = 4125600
$ws.Range("D81").Value = 1158000
$ws.Range("E81").Value = 1227600
$ws.Range("F81").Value = 1288000
$ws.Range("G81").Value = 1251200
$ws.Range("H81").Value = 1127400
$ws.Range("I81").Value = 959600
$ws.Range("J81").Value = 805900
$ws.Range("D83").Value = 305500
$ws.Range("E83").Value = 334100
$ws.Range("F83").Value = 216400
$ws.Range("G83").Value = 185400
$ws.Range("H83").Value = 160600
$ws.Range("I83").Value = 156700
$ws.Range("J83").Value = 146500
$ws.Range("D89").Value = 1218000
$ws.Range("E89").Value = 1341500
$ws.Range("F89").Value = 1140500
$ws.Range("G89").Value = 1131700
$ws.Range("H89").Value = 981800
$ws.Range("I89").Value = 1018300
$ws.Range("J89").Value = 579500
$ws.Range("D91").Value = -316200
$ws.Range("E91").Value = -301500
$ws.Range("F91").Value = -201700
$ws.Range("G91").Value = -183100
$ws.Range("H91").Value = -128900
$ws.Range("I91").Value = -153500
$ws.Range("J91").Value = -187600
$ws.Range("D94").Value = 514500
$ws.Range("E94").Value = -1681500
$ws.Range("F94").Value = -1997700
$ws.Range("G94").Value = -369800
$ws.Range("H94").Value = -40100
$ws.Range("I94").Value = -832500
$ws.Range("J94").Value = -116500
$ws.Range("D96").Value = -78400
$ws.Range("E96").Value = -126300
$ws.Range("F96").Value = -513200
$ws.Range("G96").Value = -426400
$ws.Range("H96").Value = -336500
$ws.Range("I96").Value = -247000
$ws.Range("J96").Value = -249100
$ws.Range("D100").Value = -1879500
$ws.Range("E100").Value = -329000
$ws.Range("F100").Value = -22900
$ws.Range("G100").Value = -126800
$ws.Range("H100").Value = -505700
$ws.Range("I100").Value = -97200
$ws.Range("J100").Value = -251600
$ws.Range("D101").Value = 5400
$ws.Range("E101").Value = -20400
$ws.Range("F101").Value = 7900
$ws.Range("G101").Value = 8500
$ws.Range("I101").Value = 11400
$ws.Range("J101").Value = 24300
$ws.Range("D102").Value = -141600
$ws.Range("E102").Value = -689400
$ws.Range("F102").Value = -872200
$ws.Range("G102").Value = 643600
$ws.Range("H102").Value = 435000
$ws.Range("I102").Value = 100000
$ws.Range("J102").Value = 235700
